$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.841467
$ws.Range("H2").Value = 26.524401
$ws.Range("I2").Value = 0.5917001192060068
$ws.Range("J2").Value = 0.5917001192060067
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.572065666666667
$ws.Range("N2").Value = 4.716197
$ws.Range("O2").Value = 0.1759712293834306
$ws.Range("P2").Value = 0.1759712293834305
$ws.Range("Q2").Value = 13.89936671366633
$ws.Range("R2").Value = 125.094300422997
$ws.Range("S2").Value = 0.1041221974030034
$ws.Range("T2").Value = 0.1041221974030034

# Row 3
$ws.Range("G3").Value = 8.841467
$ws.Range("H3").Value = 26.524401
$ws.Range("I3").Value = 0.5917001192060068
$ws.Range("J3").Value = 0.5917001192060067
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.237801
$ws.Range("N3").Value = 12.713403
$ws.Range("O3").Value = 0.4743638053196239
$ws.Range("P3").Value = 0.4743638053196239
$ws.Range("Q3").Value = 37.468377694067
$ws.Range("R3").Value = 337.2153992466031
$ws.Range("S3").Value = 0.2806811201546365
$ws.Range("T3").Value = 0.2806811201546364

# Row 4
$ws.Range("G4").Value = 8.841467
$ws.Range("H4").Value = 26.524401
$ws.Range("I4").Value = 0.5917001192060068
$ws.Range("J4").Value = 0.5917001192060067
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.123785
$ws.Range("N4").Value = 9.371354999999999
$ws.Range("O4").Value = 0.3496649652969456
$ws.Range("P4").Value = 0.3496649652969455
$ws.Range("Q4").Value = 27.618841992595
$ws.Range("R4").Value = 248.569577933355
$ws.Range("S4").Value = 0.2068968016483669
$ws.Range("T4").Value = 0.2068968016483669

# Row 5
$ws.Range("G5").Value = 4.103438
$ws.Range("H5").Value = 12.310314
$ws.Range("I5").Value = 0.2746155987184545
$ws.Range("J5").Value = 0.2746155987184545
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.572065666666667
$ws.Range("N5").Value = 4.716197
$ws.Range("O5").Value = 0.1759712293834306
$ws.Range("P5").Value = 0.1759712293834305
$ws.Range("Q5").Value = 6.450873995095334
$ws.Range("R5").Value = 58.05786595585799
$ws.Range("S5").Value = 0.04832444451435328
$ws.Range("T5").Value = 0.04832444451435327

# Row 6
$ws.Range("G6").Value = 4.103438
$ws.Range("H6").Value = 12.310314
$ws.Range("I6").Value = 0.2746155987184545
$ws.Range("J6").Value = 0.2746155987184545
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.237801
$ws.Range("N6").Value = 12.713403
$ws.Range("O6").Value = 0.4743638053196239
$ws.Range("P6").Value = 0.4743638053196239
$ws.Range("Q6").Value = 17.389553659838
$ws.Range("R6").Value = 156.505982938542
$ws.Range("S6").Value = 0.1302677004082129
$ws.Range("T6").Value = 0.1302677004082129

# Row 7
$ws.Range("G7").Value = 4.103438
$ws.Range("H7").Value = 12.310314
$ws.Range("I7").Value = 0.2746155987184545
$ws.Range("J7").Value = 0.2746155987184545
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.123785
$ws.Range("N7").Value = 9.371354999999999
$ws.Range("O7").Value = 0.3496649652969456
$ws.Range("P7").Value = 0.3496649652969455
$ws.Range("Q7").Value = 12.81825807283
$ws.Range("R7").Value = 115.36432265547
$ws.Range("S7").Value = 0.09602345379588834
$ws.Range("T7").Value = 0.09602345379588831

# Row 8
$ws.Range("G8").Value = 1.997574666666667
$ws.Range("H8").Value = 5.992724
$ws.Range("I8").Value = 0.1336842820755386
$ws.Range("J8").Value = 0.1336842820755386
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.572065666666667
$ws.Range("N8").Value = 4.716197
$ws.Range("O8").Value = 0.1759712293834306
$ws.Range("P8").Value = 0.1759712293834305
$ws.Range("Q8").Value = 3.140318550069778
$ws.Range("R8").Value = 28.262866950628
$ws.Range("S8").Value = 0.02352458746607384
$ws.Range("T8").Value = 0.02352458746607384

# Row 9
$ws.Range("G9").Value = 1.997574666666667
$ws.Range("H9").Value = 5.992724
$ws.Range("I9").Value = 0.1336842820755386
$ws.Range("J9").Value = 0.1336842820755386
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.237801
$ws.Range("N9").Value = 12.713403
$ws.Range("O9").Value = 0.4743638053196239
$ws.Range("P9").Value = 0.4743638053196239
$ws.Range("Q9").Value = 8.465323919974667
$ws.Range("R9").Value = 76.18791527977201
$ws.Range("S9").Value = 0.06341498475677448
$ws.Range("T9").Value = 0.06341498475677447

# Row 10
$ws.Range("G10").Value = 1.997574666666667
$ws.Range("H10").Value = 5.992724
$ws.Range("I10").Value = 0.1336842820755386
$ws.Range("J10").Value = 0.1336842820755386
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.123785
$ws.Range("N10").Value = 9.371354999999999
$ws.Range("O10").Value = 0.3496649652969456
$ws.Range("P10").Value = 0.3496649652969455
$ws.Range("Q10").Value = 6.239993780113333
$ws.Range("R10").Value = 56.15994402102
$ws.Range("S10").Value = 0.0467447098526903
$ws.Range("T10").Value = 0.04674470985269027

